$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting the existing "Tipo"/"single" column to E
$ws.Columns.Item(4).Insert()

# Copy header formatting (border/font/alignment) from the neighboring header cell (C1) to the new D1
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)

# Set the new header and value
$ws.Cells.Item(1, 4).Value = "MAE"
$ws.Cells.Item(2, 4).Value = 0.321570631223101
